# Update the AMS [N] measurements (column C) on the "Ark1" sheet to reflect
# the improved muscle geometry (internal/external obliques and latissimus
# dorsi discretization). Column D holds =Cx/C5 formulas, so their cached
# values (and the chart's numCache of 'Ark1'!$D$2:$D$9) will recompute
# automatically once the underlying C values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("C2").Value = 99.767129999999995
$ws.Range("C3").Value = 267.94510000000002
$ws.Range("C4").Value = 457.38409999999999
$ws.Range("C5").Value = 532.86289999999997
$ws.Range("C6").Value = 1125.6400000000001
$ws.Range("C7").Value = 1097.0319999999999
$ws.Range("C8").Value = 2338.2159999999999
$ws.Range("C9").Value = 1865.671

# Move the active selection from H7 to C10, matching the saved cursor
# position recorded in the workbook after the edit.
$ws.Activate()
$ws.Range("C10").Select()
